$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5
    3  = -2
    4  = 10
    5  = -8
    6  = -3
    8  = -3
    9  = 3
    11 = -1
    12 = -6
    13 = 4
    14 = -8
    15 = 2
    16 = 9
    17 = 2
    18 = -5
    19 = -4
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = -6
    25 = -2
    26 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
